# Remove the trailing "Ver no Jupiter..." / copyright boilerplate block that
# the site generator used to append after the "Requisitos" section, along
# with the blank paragraph that separated it from the requirement text.
#
# Before:
#   ...
#   LOQ4002: Reatores Quimicos (Requisito fraco)
#   <blank>
#   Ver no Jupiter Salvar em pdf Salvar em docx
#   (c) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution
#   <blank>
#   <page break>
#
# After:
#   ...
#   LOQ4002: Reatores Quimicos (Requisito fraco)
#   <blank>
#   <page break>

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter..." paragraph by its text (index lookups are
# more reliable here than the .Previous/.Next paragraph properties).
$marker = "Ver no Jupiter Salvar em pdf Salvar em docx"
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq ($marker + "`r")) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 1) {
    # The blank paragraph right before "Ver no Jupiter..." and the copyright
    # paragraph right after it are the other two paragraphs to remove.
    $blankBefore = $d.Paragraphs.Item($targetIndex - 1)
    $afterCopyright = $d.Paragraphs.Item($targetIndex + 2)

    # Delete everything from the start of the blank paragraph up to (but not
    # including) the paragraph that follows the copyright notice. That
    # consumes the blank paragraph, the "Ver no Jupiter..." paragraph, the
    # copyright paragraph, and their paragraph marks in one shot, leaving the
    # remaining blank paragraph / page-break paragraph untouched.
    $deleteRange = $d.Range($blankBefore.Range.Start, $afterCopyright.Range.Start)
    $deleteRange.Delete()
}
